# Update cryptos list values (price + 1h volume change) per the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.860.13"
$ws.Range("E2").Value = "  +0.85%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.571.71"
$ws.Range("E3").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "564.88"
$ws.Range("E5").Value = "  +3.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.58"
$ws.Range("E6").Value = "  -0.46%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.595"
$ws.Range("E8").Value = "  +1.96%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.576.76"
$ws.Range("E9").Value = "  -0.22%  "
$ws.Range("E10").Value = "  -1.68%  "
$ws.Range("E11").Value = "  +3.36%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.151"
$ws.Range("E12").Value = "  +8.04%  "
$ws.Range("E13").Value = "  +3.15%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.024.90"
$ws.Range("E14").Value = "  -0.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "58.940.54"
$ws.Range("E15").Value = "  +1.09%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.96"
$ws.Range("E16").Value = "  +7.12%  "
$ws.Range("E17").Value = "  +5.16%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.573.97"
$ws.Range("E18").Value = "  -0.62%  "
$ws.Range("E19").Value = "  +1.84%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "334.87"
$ws.Range("E20").Value = "  +0.59%  "
$ws.Range("E21").Value = "  +1.65%  "
$ws.Range("E22").Value = "  +1.71%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.59"
$ws.Range("E24").Value = "  -4.29%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.446"
$ws.Range("E25").Value = "  +6.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.161"
$ws.Range("E27").Value = "  +2.58%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.22"
$ws.Range("E28").Value = "  +2.79%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0780"
$ws.Range("E29").Value = "  +6.50%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("E31").Value = "  +1.83%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "159.10"
$ws.Range("E32").Value = "  +3.02%  "
$ws.Range("E33").Value = "  +3.03%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.91"
$ws.Range("E34").Value = "  +0.64%  "
$ws.Range("E35").Value = "  +2.74%  "
$ws.Range("E36").Value = "  +4.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.878"
$ws.Range("E37").Value = "  +8.21%  "
$ws.Range("E38").Value = "  +3.77%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.80"
$ws.Range("E39").Value = "  -0.83%  "
$ws.Range("E40").Value = "  +4.70%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "290.55"
$ws.Range("E41").Value = "  +4.47%  "
$ws.Range("E42").Value = "  +1.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0969"
$ws.Range("E44").Value = "  +3.40%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.594"
$ws.Range("E45").Value = "  +0.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.01"
$ws.Range("E48").Value = "  +3.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "124.69"
$ws.Range("E49").Value = "  +12.13%  "
$ws.Range("E50").Value = "  +1.91%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.51"
$ws.Range("E51").Value = "  +4.91%  "

# Row 46 ("Hedera") and Row 47 ("WhiteBITCoin") swapped ranking positions;
# overwrite B/C/D/E for both rows with their new (swapped) content.
$ws.Range("B46").Value = "WhiteBITCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.60"
$ws.Range("E46").Value = "  -0.36%  "

$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0533"
$ws.Range("E47").Value = "  +1.49%  "
